$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 and Row 10 had their Id (A), Antal (I), Ost (Q) and Nord (R) values
# swapped between each other. Read current values first (Value2 avoids the
# text-typed Antal column being mis-detected), then write the swapped
# values back, preserving the text type of column I with a leading
# apostrophe (Excel's "store as text" convention) since it holds numeric-
# looking text rather than a real number.

$a9 = $ws.Range("A9").Value2
$i9 = $ws.Range("I9").Text
$q9 = $ws.Range("Q9").Value2
$r9 = $ws.Range("R9").Value2

$a10 = $ws.Range("A10").Value2
$i10 = $ws.Range("I10").Text
$q10 = $ws.Range("Q10").Value2
$r10 = $ws.Range("R10").Value2

$ws.Range("A9").Value = $a10
$ws.Range("I9").Value = "'" + $i10
$ws.Range("Q9").Value = $q10
$ws.Range("R9").Value = $r10

$ws.Range("A10").Value = $a9
$ws.Range("I10").Value = "'" + $i9
$ws.Range("Q10").Value = $q9
$ws.Range("R10").Value = $r9
